$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.820.67"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "3.113.56"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "3.110.86"
$ws.Range("E8").Value = "  +1.43%  "

$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "

$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "3.631.37"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "66.831.40"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").Value = "3.115.37"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "475.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.713"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("E26").Value = "  -1.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.61%  "

$ws.Range("E34").Value = "  -7.61%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.975"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "386.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.32%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.829.29"
$ws.Range("E45").Value = "  +2.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0356"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("E51").Value = "  -1.30%  "

